# Generate Report for handoff
#
# The "acb8a402-b60a-4eac-b7c6-6347f3b08d95" row (row 3 on every sheet) moves
# from "Handed back: in sync with en-US" to "Ready for handoff", and its
# "Latest Handoff Datetime" is refreshed with the new handoff timestamp.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: Status columns for the acb8a402 row (row 3)
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# zh-cn detail sheet: Status + Latest Handoff Datetime for the acb8a402 row
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-02-15 03:46:15"

# de-de detail sheet: Status + Latest Handoff Datetime for the acb8a402 row
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-02-15 03:46:32"
